$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E values (rows 2-51).
# Rows 2-4 originally carried style index 2; the edited values drop that
# explicit style so the cell reverts to the default ("Normal") style.
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").Value = 2.23

$ws.Range("E3").Style = "Normal"
$ws.Range("E3").Value = 2.92

$ws.Range("E4").Style = "Normal"
$ws.Range("E4").Value = 9.6300000000000008

$ws.Range("E5").Value = 11.56
$ws.Range("E6").Value = 3.29
$ws.Range("E7").Value = 1.47
$ws.Range("E8").Value = 6.28
$ws.Range("E9").Value = 11.98
$ws.Range("E10").Value = 5.62
$ws.Range("E11").Value = 7.41
$ws.Range("E12").Value = 5.0999999999999996
$ws.Range("E13").Value = 5.15
$ws.Range("E14").Value = 5.66
$ws.Range("E15").Value = 2.4700000000000002
$ws.Range("E16").Value = 4.07
$ws.Range("E17").Value = 5.29
$ws.Range("E18").Value = 5.0999999999999996
$ws.Range("E19").Value = 1.77
$ws.Range("E20").Value = 1.74
$ws.Range("E21").Value = 4.21
$ws.Range("E22").Value = 26.21
$ws.Range("E23").Value = 0.86
$ws.Range("E24").Value = 9.51
$ws.Range("E25").Value = 4.6500000000000004
$ws.Range("E26").Value = 14.23
$ws.Range("E27").Value = 2.87
$ws.Range("E28").Value = 5.18
$ws.Range("E29").Value = 3.21
$ws.Range("E30").Value = 3.02
$ws.Range("E31").Value = 12.25
$ws.Range("E32").Value = 2.06
$ws.Range("E33").Value = 2.34
$ws.Range("E34").Value = 6.42
$ws.Range("E35").Value = 52.28
$ws.Range("E36").Value = 1.66
$ws.Range("E37").Value = 2.15
$ws.Range("E38").Value = 34.97
$ws.Range("E39").Value = 1.92
$ws.Range("E40").Value = 10.37
$ws.Range("E41").Value = 9.52
$ws.Range("E42").Value = 7.24
$ws.Range("E43").Value = 1.22
$ws.Range("E44").Value = 5.22
$ws.Range("E45").Value = 3.83
$ws.Range("E46").Value = 2.09
$ws.Range("E47").Value = 5.2
$ws.Range("E48").Value = 12.38
$ws.Range("E49").Value = 2.96
$ws.Range("E50").Value = 4.67
$ws.Range("E51").Value = 0.98

# The author left the selection on A11 when saving.
$ws.Range("A11").Select()
